$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list on Fri Jul  7 19:31:00 UTC 2023 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) for rows 2-51.
# Price values that Excel would otherwise auto-convert to numbers are
# forced back to plain text first, matching the original text-formatted data.

$ws.Range("D2").Value = "30.208.96"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "1.862.59"
$ws.Range("E3").Value = "  -1.44%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.86"
$ws.Range("E5").Value = "  -1.46%  "
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4669"
$ws.Range("E7").Value = "  -0.56%  "
$ws.Range("D8").Value = "0.2828"
$ws.Range("E8").Value = "  -0.76%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06532"
$ws.Range("E9").Value = "  -1.33%  "
$ws.Range("D10").Value = "21.34"
$ws.Range("E10").Value = "  +2.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07850"
$ws.Range("E11").Value = "  +0.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.32"
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("D13").Value = "1.865.65"
$ws.Range("E13").Value = "  -1.20%  "
$ws.Range("D14").Value = "5.093"
$ws.Range("E14").Value = "  -0.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6718"
$ws.Range("E15").Value = "  -0.72%  "
$ws.Range("D16").Value = "279.19"
$ws.Range("E16").Value = "  -2.12%  "
$ws.Range("D17").Value = "30.197.32"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.513"
$ws.Range("E19").Value = "  +1.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.62"
$ws.Range("E20").Value = "  -0.70%  "
$ws.Range("D21").Value = "2.115.00"
$ws.Range("E21").Value = "  -0.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.000007276"
$ws.Range("E22").Value = "  -0.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "6.149"
$ws.Range("E24").Value = "  -0.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.201"
$ws.Range("E25").Value = "  -2.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.65"
$ws.Range("E26").Value = "  -1.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.15"
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.920"
$ws.Range("E28").Value = "  -4.12%  "
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09679"
$ws.Range("E30").Value = "  -0.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.419"
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("D32").Value = "1.474"
$ws.Range("E32").Value = "  -1.08%  "
$ws.Range("D33").Value = "4.077"
$ws.Range("E33").Value = "  -2.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04689"
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("D35").Value = "1.113"
$ws.Range("E35").Value = "  +1.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7049"
$ws.Range("E36").Value = "  -1.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.728"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01853"
$ws.Range("E38").Value = "  -1.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.534"
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.219"
$ws.Range("E40").Value = "  -7.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.22"
$ws.Range("E41").Value = "  +0.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.940"
$ws.Range("E42").Value = "  -2.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8460"
$ws.Range("E43").Value = "  -2.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "103.99"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4162"
$ws.Range("E46").Value = "  -1.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.188"
$ws.Range("E47").Value = "  -1.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.190"
$ws.Range("E48").Value = "  -0.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "936.15"
$ws.Range("E49").Value = "  -7.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.05"
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("E51").Value = "  -2.82%  "
